# Re-run of the asthma IR analysis with new data (other states removed).
# Updates the computed income-sum column (C) with the refreshed results,
# narrows the selected range to the recalculated column, and forces the
# worksheet's print orientation (matching a vertical pageSetup entry that
# newer Excel adds on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated results in column C (rows 2-11) ------------------------------
$ws.Range("C2").Value  = 132643.90259411329
$ws.Range("C3").Value  = 13848.317991206073
$ws.Range("C4").Value  = 9284.5334103479945
$ws.Range("C5").Value  = 109511.0511925592
$ws.Range("C6").Value  = 5740.3580004620526
$ws.Range("C7").Value  = 24807.243343738122
$ws.Range("C8").Value  = 32352.087319831586
$ws.Range("C9").Value  = 37595.732117752814
$ws.Range("C10").Value = 32120.997452340147
$ws.Range("C11").Value = 27.484359988560563

# --- Selection now only spans the recalculated column C -------------------
[void]$ws.Range("C2:C10").Select()

# --- Page setup: force portrait orientation (adds <pageSetup .../>) -------
$ws.PageSetup.Orientation = 1
